# Apply updated PoS tag labels to the "Arabic", "Hebrew" and "Algerian"
# tag-stats tables (rows 18-21 of the TRAIN/DEV/TEST/TOTAL summary blocks)
# after re-running the Indonesian PoS and PUD tests.

$wb = $excel.ActiveWorkbook

# --- Arabic sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Arabic")
$ws.Range("K20").Value = "O"
$ws.Range("K21").Value = "INTJ"

# --- Hebrew sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Hebrew")
$ws.Range("F18").Value = "X"
$ws.Range("K18").Value = "PART"
$ws.Range("A19").Value = "O"
$ws.Range("F19").Value = "PART"
$ws.Range("P19").Value = "PART"
$ws.Range("A20").Value = "PART"
$ws.Range("F20").Value = "SYM"
$ws.Range("P20").Value = "SYM"
$ws.Range("A21").Value = "SYM"
$ws.Range("F21").Value = "O"
$ws.Range("K21").Value = "INTJ"
$ws.Range("P21").Value = "O"

# --- Algerian sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("Algerian")
$ws.Range("K18").Value = "NUM"
$ws.Range("F19").Value = "O"
$ws.Range("K19").Value = "SYM"
$ws.Range("A20").Value = "O"
$ws.Range("K20").Value = "AUX"
$ws.Range("P20").Value = "O"
$ws.Range("A21").Value = "SYM"
$ws.Range("F21").Value = "NUM"
$ws.Range("K21").Value = "O"
$ws.Range("P21").Value = "SYM"
